# Apply updated NATMI LR-pair statistics (Matn1-Itgb1) per Dr Hou's advice.
# Updates columns E, G, H, K, M, N, O, P, Q, R, S, T for data rows 2-5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{
        E = 3
        G = 0.3987123333333333
        H = 1.196137
        K = 3
        M = 98.946724
        N = 296.840172
        O = 0.2098009692989996
        P = 0.2098009692989996
        Q = 39.45127920172933
        R = 355.061512815564
        S = 0.2098009692989996
        T = 0.2098009692989996
    }
    3 = @{
        E = 3
        G = 0.3987123333333333
        H = 1.196137
        K = 3
        M = 163.0062356666667
        N = 489.018707
        O = 0.345629090707923
        P = 0.3456290907079231
        Q = 64.99259657053989
        R = 584.933369134859
        S = 0.345629090707923
        T = 0.3456290907079231
    }
    4 = @{
        E = 3
        G = 0.3987123333333333
        H = 1.196137
        K = 3
        M = 65.39610666666668
        N = 196.18832
        O = 0.1386621609326595
        P = 0.1386621609326595
        Q = 26.07423427998223
        R = 234.66810851984
        S = 0.1386621609326595
        T = 0.1386621609326595
    }
    5 = @{
        E = 3
        G = 0.3987123333333333
        H = 1.196137
        K = 3
        M = 144.2727966666667
        N = 432.81839
        O = 0.3059077790604178
        P = 0.3059077790604179
        Q = 57.52334339549222
        R = 517.71009055943
        S = 0.3059077790604178
        T = 0.3059077790604179
    }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
